# Update TPM-derived values in the LR-pairs sheet (Vip-Vipr2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3882076666666667
$ws.Range("H2").Value = 1.164623
$ws.Range("M2").Value = 0.01962266666666667
$ws.Range("N2").Value = 0.058868
$ws.Range("O2").Value = 0.02207703711370904
$ws.Range("P2").Value = 0.02207703711370904
$ws.Range("Q2").Value = 0.007617669640444445
$ws.Range("R2").Value = 0.068559026764
$ws.Range("S2").Value = 0.02207703711370904
$ws.Range("T2").Value = 0.02207703711370904

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3882076666666667
$ws.Range("H3").Value = 1.164623
$ws.Range("O3").Value = 0.7491195324474467
$ws.Range("P3").Value = 0.7491195324474467
$ws.Range("Q3").Value = 0.2584832869554444
$ws.Range("R3").Value = 2.326349582599
$ws.Range("S3").Value = 0.7491195324474467
$ws.Range("T3").Value = 0.7491195324474467

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3882076666666667
$ws.Range("H4").Value = 1.164623
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1131433333333333
$ws.Range("N4").Value = 0.33943
$ws.Range("O4").Value = 0.1272951129222372
$ws.Range("P4").Value = 0.1272951129222372
$ws.Range("Q4").Value = 0.04392310943222223
$ws.Range("R4").Value = 0.39530798489
$ws.Range("S4").Value = 0.1272951129222372
$ws.Range("T4").Value = 0.1272951129222372

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3882076666666667
$ws.Range("H5").Value = 1.164623
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09022333333333332
$ws.Range("N5").Value = 0.27067
$ws.Range("O5").Value = 0.1015083175166071
$ws.Range("P5").Value = 0.1015083175166071
$ws.Range("Q5").Value = 0.03502538971222222
$ws.Range("R5").Value = 0.31522850741
$ws.Range("S5").Value = 0.1015083175166071
$ws.Range("T5").Value = 0.1015083175166071
